$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell, forcing Text number format first
# whenever the new value would otherwise be auto-parsed by Excel as a
# number (which would corrupt literal representations such as
# '212.10', '4.00', '26.169.43', or multi-dot strings).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range("D2").Value = '26.169.43'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '1.606.47'
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("E4").Value = '  -0.14%  '

Set-TextValue $ws.Range("D5") '212.10'
$ws.Range("E5").Value = '  +1.38%  '

Set-TextValue $ws.Range("D7") '0.481'
$ws.Range("E7").Value = '  +0.33%  '

Set-TextValue $ws.Range("D8") '0.249'
$ws.Range("E8").Value = '  +1.17%  '

Set-TextValue $ws.Range("D9") '0.0617'
$ws.Range("E9").Value = '  +1.31%  '

Set-TextValue $ws.Range("D10") '18.11'
$ws.Range("E10").Value = '  +1.28%  '

Set-TextValue $ws.Range("D11") '0.0795'
$ws.Range("E11").Value = '  +1.39%  '

$ws.Range("D12").Value = '1.831.75'
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").Value = '1.608.32'
$ws.Range("E13").Value = '  +0.47%  '

Set-TextValue $ws.Range("D14") '4.00'
$ws.Range("E14").Value = '  -1.07%  '

Set-TextValue $ws.Range("D15") '0.509'
$ws.Range("E15").Value = '  -0.05%  '

$ws.Range("D16").Value = '26.169.73'
$ws.Range("E16").Value = '  +1.34%  '

Set-TextValue $ws.Range("D17") '60.60'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  +1.56%  '

$ws.Range("E19").Value = '  -0.19%  '

Set-TextValue $ws.Range("D20") '199.24'
$ws.Range("E20").Value = '  +5.13%  '

Set-TextValue $ws.Range("D21") '4.24'
$ws.Range("E21").Value = '  +1.39%  '

Set-TextValue $ws.Range("D22") '9.41'
$ws.Range("E22").Value = '  +0.74%  '

Set-TextValue $ws.Range("D23") '6.00'
$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("E24").Value = '  +3.21%  '

Set-TextValue $ws.Range("D25") '141.88'
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").Value = '  +2.08%  '

$ws.Range("E27").Value = '  -0.24%  '

Set-TextValue $ws.Range("D28") '15.17'
$ws.Range("E28").Value = '  +1.48%  '

Set-TextValue $ws.Range("D29") '6.48'
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("E30").Value = '  -1.12%  '

Set-TextValue $ws.Range("D31") '0.0473'
$ws.Range("E31").Value = '  +1.21%  '

Set-TextValue $ws.Range("D32") '3.13'
$ws.Range("E32").Value = '  +1.82%  '

Set-TextValue $ws.Range("D33") '3.02'
$ws.Range("E33").Value = '  +0.68%  '

$ws.Range("E34").Value = '  +2.22%  '

Set-TextValue $ws.Range("D35") '2.36'
$ws.Range("E35").Value = '  -1.64%  '

$ws.Range("D36").Value = '1.107.98'
$ws.Range("E36").Value = '  +0.99%  '

Set-TextValue $ws.Range("D37") '2.37'
$ws.Range("E37").Value = '  -0.01%  '

Set-TextValue $ws.Range("D40") '0.502'
$ws.Range("E40").Value = '  +0.78%  '

$ws.Range("E41").Value = '  -0.69%  '

Set-TextValue $ws.Range("D42") '0.784'
$ws.Range("E42").Value = '  +5.68%  '

$ws.Range("D43").Value = '1.744.84'
$ws.Range("E43").Value = '  +0.55%  '

Set-TextValue $ws.Range("D44") '5.14'
$ws.Range("E44").Value = '  +1.43%  '

Set-TextValue $ws.Range("D45") '92.98'
$ws.Range("E45").Value = '  -2.82%  '

Set-TextValue $ws.Range("D48") '53.58'
$ws.Range("E48").Value = '  +0.73%  '

Set-TextValue $ws.Range("D49") '0.0509'
$ws.Range("E49").Value = '  -0.40%  '

$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("E51").Value = '  -0.08%  '

# Row swaps: VeChain/PaxDollar (rows 38-39) and BabyDogeCoin/RenderToken (rows 46-47)
$ws.Range("B38").Value = 'PaxDollar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D38") '1.00'
$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.0152'
$ws.Range("E39").Value = '  +0.59%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D46") '1.54'
$ws.Range("E46").Value = '  +7.64%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0103'
$ws.Range("E47").Value = '  -7.89%  '

